# Update "Generate Report for Handback" timestamps across the workbook.
# Shared string "2016-08-26 17:06:11" is used by both Overview!G2 (Latest HO
# Xliff Generate Date) and de-de!H2 (Correspond Handoff Datetime), so both
# cells must be updated together to the new timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview.Range("G2").Value = "2016-08-26 17:06:57"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime.
$wsZhCn.Range("H2").Value = "2016-08-26 17:06:53"
$wsZhCn.Range("K2").Value = "2016-08-26 17:07:16"

# de-de sheet: Correspond Handoff Datetime (shared value with Overview!G2)
# and Correspond Handback DateTime.
$wsDeDe.Range("H2").Value = "2016-08-26 17:06:57"
$wsDeDe.Range("K2").Value = "2016-08-26 17:07:23"
